$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")

# ---------------------------------------------------------------------------
# 1) Structural change: the two body-fat rows (FR_android / FR_gynoid) are
#    consolidated into a single row (FG_Prozent). Deleting one of the pair
#    shifts everything below up by one row, which also removes the last
#    (empty) row of the sheet (94 -> 93), matching the new dimension.
# ---------------------------------------------------------------------------
$ws.Rows.Item(9).EntireRow.Delete()

# ---------------------------------------------------------------------------
# 2) Update cell text so the table reflects the revised data dictionary.
# ---------------------------------------------------------------------------

# Row 7: waist circumference -> German label "Taillenumfang"
$ws.Cells.Item(7,3).Value2 = "Taillenumfang"

# Row 8: replace the old %body-fat (android) row with the new body-fat % row
$ws.Cells.Item(8,2).Value2 = "FG_Prozent"
$ws.Cells.Item(8,3).Value2 = "Gewebe%Fett " + [char]0x2013 + " Gesamt in %"

# Rows 9-10 (FMI, Alter_BE) keep their text - unchanged.

# Row 20: the "Total sugar intake" row becomes the "Sugars_NCI" row (name and
# label now hold the same variable-style text).
$ws.Cells.Item(20,2).Value2 = "Sugars_NCI"
$ws.Cells.Item(20,3).Value2 = "Sugars_NCI"

# Rows 11-19: shorten / simplify the nutrient intake labels
$ws.Cells.Item(11,3).Value2 = "Energy " + [char]0x2013 + " incl. energy from dietary fiber"
$ws.Cells.Item(12,3).Value2 = "Carbohydrates"
$ws.Cells.Item(13,3).Value2 = "Protein"
$ws.Cells.Item(14,3).Value2 = "Fat"
$ws.Cells.Item(15,3).Value2 = "Alcohol"
$ws.Cells.Item(16,3).Value2 = "Dietary Fiber"
$ws.Cells.Item(17,3).Value2 = "Saturated Fat"
$ws.Cells.Item(18,3).Value2 = "Monounsaturated Fat"
$ws.Cells.Item(19,3).Value2 = "Polyunsaturated Fat"

# Rows 21-24: shorten remaining nutrient labels
$ws.Cells.Item(21,3).Value2 = "Glucose"
$ws.Cells.Item(22,3).Value2 = "Fructose"
$ws.Cells.Item(23,3).Value2 = "Sodium"
$ws.Cells.Item(24,3).Value2 = "Potassium"

# Rows 19-24 now also carry a "decimal" valueType in column D (previously
# blank for this block of rows).
for ($r = 19; $r -le 24; $r++) {
  $ws.Cells.Item($r,4).Value2 = "decimal"
}

# ---------------------------------------------------------------------------
# 3) Formatting touch-ups that accompanied the content rewrite.
# ---------------------------------------------------------------------------

# Header row becomes bold.
$ws.Range("A1:D1").Font.Bold = $true

# The rewritten label column (rows 7-8 and 11-22) is left/top aligned with
# wrap text, consistent with how the new text was pasted in.
$ws.Range("C7:C8").WrapText = $true
$ws.Range("C7:C8").HorizontalAlignment = -4131
$ws.Range("C7:C8").VerticalAlignment = -4160

$ws.Range("C11:C19").WrapText = $true
$ws.Range("C11:C19").HorizontalAlignment = -4131
$ws.Range("C11:C19").VerticalAlignment = -4160

$ws.Range("B20:C20").WrapText = $true

$ws.Range("C21:C22").WrapText = $true

# ---------------------------------------------------------------------------
# 4) View state: zoomed out a bit, with a different active selection.
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 70
$ws.Range("B28").Select()
